$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (column D) value (or $null if unchanged),
# whether the Price value must be forced to text (to avoid Excel's numeric
# auto-conversion stripping formatting such as trailing zeros), and the new
# Volume(1h) (column E) value (column E always changes).
$updates = @(
    @{Row=2; D='23.959.11'; ForceText=$false; E='  +0.67%  '}
    @{Row=3; D='1.652.82'; ForceText=$false; E='  +2.24%  '}
    @{Row=4; D='1.002'; ForceText=$true; E='  -0.06%  '}
    @{Row=5; D='309.38'; ForceText=$true; E='  +0.47%  '}
    @{Row=6; D='1.000'; ForceText=$true; E='  -0.17%  '}
    @{Row=7; D='0.3903'; ForceText=$true; E='  -0.74%  '}
    @{Row=8; D='0.3840'; ForceText=$true; E='  +0.24%  '}
    @{Row=9; D='51.26'; ForceText=$true; E='  +3.91%  '}
    @{Row=10; D='1.357'; ForceText=$true; E='  +0.42%  '}
    @{Row=11; D=$null; ForceText=$false; E='  -0.07%  '}
    @{Row=12; D='0.08452'; ForceText=$true; E='  +0.26%  '}
    @{Row=13; D='24.01'; ForceText=$true; E='  +1.62%  '}
    @{Row=14; D='7.123'; ForceText=$true; E='  +1.36%  '}
    @{Row=15; D='7.875'; ForceText=$true; E='  +4.06%  '}
    @{Row=16; D='0.00001317'; ForceText=$true; E='  +3.18%  '}
    @{Row=17; D='1.655.21'; ForceText=$false; E='  +2.27%  '}
    @{Row=18; D='94.65'; ForceText=$true; E='  +0.91%  '}
    @{Row=19; D='0.06988'; ForceText=$true; E='  +0.84%  '}
    @{Row=20; D='19.81'; ForceText=$true; E='  -0.65%  '}
    @{Row=21; D='6.932'; ForceText=$true; E='  +1.92%  '}
    @{Row=22; D='1.000'; ForceText=$true; E='  -0.18%  '}
    @{Row=23; D='13.65'; ForceText=$true; E='  +1.88%  '}
    @{Row=24; D='23.957.45'; ForceText=$false; E='  +0.67%  '}
    @{Row=25; D=$null; ForceText=$false; E='  +1.72%  '}
    @{Row=26; D='3.018'; ForceText=$true; E='  +6.61%  '}
    @{Row=27; D='22.13'; ForceText=$true; E='  -0.25%  '}
    @{Row=28; D='151.12'; ForceText=$true; E='  -3.73%  '}
    @{Row=29; D='5.451'; ForceText=$true; E='  +3.05%  '}
    @{Row=30; D='139.50'; ForceText=$true; E='  -0.19%  '}
    @{Row=31; D='7.881'; ForceText=$true; E='  +1.06%  '}
    @{Row=32; D='2.487'; ForceText=$true; E='  -0.13%  '}
    @{Row=33; D='1.837.80'; ForceText=$false; E='  +2.48%  '}
    @{Row=34; D='1.045'; ForceText=$true; E='  +7.11%  '}
    @{Row=35; D='0.08115'; ForceText=$true; E='  +0.32%  '}
    @{Row=36; D='0.02972'; ForceText=$true; E='  +3.34%  '}
    @{Row=37; D='6.762'; ForceText=$true; E='  +3.24%  '}
    @{Row=38; D='10.88'; ForceText=$true; E='  +5.46%  '}
    @{Row=39; D=$null; ForceText=$false; E='  +1.06%  '}
    @{Row=40; D='0.09145'; ForceText=$true; E='  +0.29%  '}
    @{Row=41; D='0.7564'; ForceText=$true; E='  +1.04%  '}
    @{Row=42; D='13.46'; ForceText=$true; E='  -0.55%  '}
    @{Row=43; D='1.427'; ForceText=$true; E='  +0.15%  '}
    @{Row=44; D='16.31'; ForceText=$true; E='  +1.75%  '}
    @{Row=45; D='0.6959'; ForceText=$true; E='  +0.90%  '}
    @{Row=46; D='2.458'; ForceText=$true; E='  -0.35%  '}
    @{Row=47; D='4.086'; ForceText=$true; E='  +0.51%  '}
    @{Row=48; D='0.9996'; ForceText=$true; E='  -0.16%  '}
    @{Row=49; D='0.08296'; ForceText=$true; E='  +0.90%  '}
    @{Row=50; D='134.66'; ForceText=$true; E='  +0.13%  '}
    @{Row=51; D='1.205'; ForceText=$true; E='  +0.30%  '}
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        if ($u.ForceText) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

"Updated $($updates.Count) rows"
